$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 1
$ws.Range("A1").Value = " "

# Row 2 - Title "GOAL SEEK"
$ws.Range("C2").Value = "GOAL SEEK"
$ws.Range("C2:G2").Merge()
$ws.Range("C2:G2").HorizontalAlignment = -4108
$ws.Range("C2").Font.Size = 18
$ws.Range("C2").Font.ThemeColor = 3
$ws.Rows("2").RowHeight = 23.25

# Row 4 - headers
$ws.Range("C4").Value = "SL NO"
$ws.Range("D4").Value = "NAME"
$ws.Range("E4").Value = "SUBJECT"
$ws.Range("F4").Value = "MARKS"

# Row 5
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = "DEEP"
$ws.Range("E5").Value = "MATH"
$ws.Range("F5").Value = 87

# Row 6
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = "DHURVI"
$ws.Range("E6").Value = "ENG"
$ws.Range("F6").Value = 87

# Row 7
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = "ITTY"
$ws.Range("E7").Value = "APTITUDE"
$ws.Range("F7").Value = 87

# Row 9 - Average
$ws.Range("E9").Value = "AVERAGE"
$ws.Range("F9").Formula = "=AVERAGE(F5:F7)"

$ws.Columns("E").ColumnWidth = 10.85546875

# Goal Seek: make F9 = 80 by changing F7
$ws.Range("F9").GoalSeek(80, $ws.Range("F7"))

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("F9").Select()
